$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = 0.2249999940395355
$ws.Range("C2").Value = 0.2156264036893845
$ws.Range("D2").Value = 0.2000000029802322
$ws.Range("E2").Value = 0.1917306184768677
$ws.Range("B3").Value = 0.2374999970197678
$ws.Range("C3").Value = 0.1901871562004089
$ws.Range("D3").Value = 0.2000000029802322
$ws.Range("E3").Value = 0.1877334713935852
$ws.Range("B4").Value = 0.3125
$ws.Range("C4").Value = 0.1881064027547836
$ws.Range("D4").Value = 0.2000000029802322
$ws.Range("E4").Value = 0.1862240731716156
$ws.Range("B5").Value = 0.3125
$ws.Range("C5").Value = 0.1852631270885468
$ws.Range("D5").Value = 0.2249999940395355
$ws.Range("E5").Value = 0.1831875294446945
$ws.Range("B6").Value = 0.3687500059604645
$ws.Range("C6").Value = 0.1830499917268753
$ws.Range("D6").Value = 0.4000000059604645
$ws.Range("E6").Value = 0.1774629205465317
$ws.Range("B7").Value = 0.4124999940395355
$ws.Range("C7").Value = 0.1787919998168945
$ws.Range("D7").Value = 0.2249999940395355
$ws.Range("E7").Value = 0.1768079698085785
$ws.Range("B8").Value = 0.449999988079071
$ws.Range("C8").Value = 0.1716659367084503
$ws.Range("D8").Value = 0.4000000059604645
$ws.Range("E8").Value = 0.1649817228317261
$ws.Range("B9").Value = 0.5062500238418579
$ws.Range("C9").Value = 0.1633847951889038
$ws.Range("D9").Value = 0.625
$ws.Range("E9").Value = 0.1425943374633789
$ws.Range("B10").Value = 0.543749988079071
$ws.Range("C10").Value = 0.1503585129976273
$ws.Range("D10").Value = 0.699999988079071
$ws.Range("E10").Value = 0.1239223033189774
$ws.Range("B11").Value = 0.6187499761581421
$ws.Range("C11").Value = 0.1373969167470932
$ws.Range("D11").Value = 0.699999988079071
$ws.Range("E11").Value = 0.1063784509897232
$ws.Range("B12").Value = 0.637499988079071
$ws.Range("C12").Value = 0.1275408267974854
$ws.Range("D12").Value = 0.699999988079071
$ws.Range("E12").Value = 0.1056616082787514
$ws.Range("B13").Value = 0.6875
$ws.Range("C13").Value = 0.1201655119657516
$ws.Range("D13").Value = 0.675000011920929
$ws.Range("E13").Value = 0.1009130254387856
$ws.Range("B14").Value = 0.675000011920929
$ws.Range("C14").Value = 0.112513855099678
$ws.Range("D14").Value = 0.699999988079071
$ws.Range("E14").Value = 0.09546820819377899
$ws.Range("B15").Value = 0.71875
$ws.Range("C15").Value = 0.1073798313736916
$ws.Range("D15").Value = 0.699999988079071
$ws.Range("E15").Value = 0.08719857037067413
$ws.Range("B16").Value = 0.762499988079071
$ws.Range("C16").Value = 0.097477987408638
$ws.Range("D16").Value = 0.7250000238418579
$ws.Range("E16").Value = 0.08891966193914413
$ws.Range("B17").Value = 0.7749999761581421
$ws.Range("C17").Value = 0.0888277143239975
$ws.Range("D17").Value = 0.824999988079071
$ws.Range("E17").Value = 0.07608769834041595
$ws.Range("B18").Value = 0.7875000238418579
$ws.Range("C18").Value = 0.08468399941921234
$ws.Range("D18").Value = 0.800000011920929
$ws.Range("E18").Value = 0.076198048889637
$ws.Range("B19").Value = 0.8187500238418579
$ws.Range("C19").Value = 0.07598400861024857
$ws.Range("D19").Value = 0.824999988079071
$ws.Range("E19").Value = 0.0726330354809761
$ws.Range("B20").Value = 0.8125
$ws.Range("C20").Value = 0.07229314744472504
$ws.Range("D20").Value = 0.8500000238418579
$ws.Range("E20").Value = 0.06865433603525162
$ws.Range("B21").Value = 0.8812500238418579
$ws.Range("C21").Value = 0.06271536648273468
$ws.Range("D21").Value = 0.7749999761581421
$ws.Range("E21").Value = 0.0709904208779335
$ws.Range("B22").Value = 0.8687499761581421
$ws.Range("C22").Value = 0.05907921120524406
$ws.Range("D22").Value = 0.8999999761581421
$ws.Range("E22").Value = 0.06305022537708282
$ws.Range("B23").Value = 0.8999999761581421
$ws.Range("C23").Value = 0.05201143026351929
$ws.Range("D23").Value = 0.8999999761581421
$ws.Range("E23").Value = 0.06148362159729004
$ws.Range("B24").Value = 0.918749988079071
$ws.Range("C24").Value = 0.04715242236852646
$ws.Range("D24").Value = 0.8999999761581421
$ws.Range("E24").Value = 0.06075760722160339
$ws.Range("B25").Value = 0.9312499761581421
$ws.Range("C25").Value = 0.04166350141167641
$ws.Range("D25").Value = 0.8500000238418579
$ws.Range("E25").Value = 0.06533880531787872
$ws.Range("B26").Value = 0.9437500238418579
$ws.Range("C26").Value = 0.03710343688726425
$ws.Range("D26").Value = 1.0
$ws.Range("E26").Value = 0.05113442987203598
$ws.Range("B27").Value = 0.949999988079071
$ws.Range("C27").Value = 0.03333047777414322
$ws.Range("D27").Value = 0.8999999761581421
$ws.Range("E27").Value = 0.05971075966954231
$ws.Range("B28").Value = 0.981249988079071
$ws.Range("C28").Value = 0.03051383793354034
$ws.Range("D28").Value = 1.0
$ws.Range("E28").Value = 0.04816677421331406
$ws.Range("B29").Value = 0.9750000238418579
$ws.Range("C29").Value = 0.02773591317236423
$ws.Range("D29").Value = 1.0
$ws.Range("E29").Value = 0.045600775629282
$ws.Range("B30").Value = 0.9750000238418579
$ws.Range("C30").Value = 0.02533054910600185
$ws.Range("D30").Value = 1.0
$ws.Range("E30").Value = 0.04405773058533669
$ws.Range("B31").Value = 0.9750000238418579
$ws.Range("C31").Value = 0.02425135299563408
$ws.Range("D31").Value = 1.0
$ws.Range("E31").Value = 0.0431615486741066
